$d = $word.ActiveDocument

# Locate the paragraph that holds the "What are values, Variables,
# Operators or function ... isn't met?" sentence (the 5th paragraph in
# this document, but search by content to be resilient).
$paraIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Operators or function*") {
        $paraIndex = $i
        break
    }
}

$full = $d.Paragraphs($paraIndex).Range
# Exclude the trailing paragraph mark so the <w:pPr> is preserved.
$target = $d.Range($full.Start, $full.End - 1)

$innerXml = '<w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>What are Operators or function</w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>s</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> that allow me to write the code I want in case the conditions for my code in   if () {}</w:t></w:r>' + `
    '<w:r w:rsidR="00807944"><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> and else </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r w:rsidR="00807944"><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>if(</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r w:rsidR="00807944"><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>) {}</w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> isn' + [char]0x2019 + 't met?</w:t></w:r>'

$wordXml = '<?xml version="1.0" standalone="yes"?>' + `
    '<?mso-application progid="Word.Document"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($wordXml)
